# Add a new bash_lib entry documenting tar's advanced / zip options
# (commit: "Add bash entry on tar x zip options").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 55

$category = "tar"
$subject  = "# Advance option"

$bodyLines = @(
    "# Zip otpion",
    "Use z in option to switch on the zip feature: ``tar -czvf foo.tar.gz bar/``",
    "> the f option must place before the file name!",
    "# Ignore certain folder",
    "``tar --exclude='.git' --exclude='target/' -czvf foo.tar.gz bar/``",
    "> Make sure the exclude declare first",
    "> This will exclude folder in that name among all levels of directory",
    "# Ignore version controls",
    "``tar --exclude-vcs ...``"
)
$body = [string]::Join("`n", $bodyLines)

$ws.Cells.Item($row, 1).Value = $category
$ws.Cells.Item($row, 2).Value = $subject
$ws.Cells.Item($row, 3).Value = $body

# Match the formatting used by the other long-form rows in column C (wrap text).
$ws.Cells.Item($row, 3).WrapText = $true

# Row height matches the rendered size of similarly long entries (e.g. row 41).
$ws.Rows.Item($row).RowHeight = 141.75

# Leave the selection where Excel lands right after typing the last cell of the row.
$ws.Range("C56").Select()
